$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "314.50"
Set-TextValue $ws.Range("E2") "2.67%"
Set-TextValue $ws.Range("E3") "-2.11%"
Set-TextValue $ws.Range("D4") "5.134"
Set-TextValue $ws.Range("E4") "0.33%"
Set-TextValue $ws.Range("E5") "2.74%"
Set-TextValue $ws.Range("D6") "2.114"
Set-TextValue $ws.Range("E6") "-0.04%"
Set-TextValue $ws.Range("D7") "4.159"
Set-TextValue $ws.Range("E7") "0.91%"
Set-TextValue $ws.Range("D8") "7.957"
Set-TextValue $ws.Range("E8") "-0.17%"
Set-TextValue $ws.Range("D9") "0.9290"
Set-TextValue $ws.Range("E9") "0.91%"
Set-TextValue $ws.Range("D10") "0.1014"
Set-TextValue $ws.Range("E10") "4.30%"
Set-TextValue $ws.Range("D11") "0.1873"
Set-TextValue $ws.Range("E11") "0.87%"
Set-TextValue $ws.Range("D12") "0.09129"
Set-TextValue $ws.Range("E12") "4.36%"
Set-TextValue $ws.Range("D13") "0.03609"
Set-TextValue $ws.Range("E13") "1.12%"
Set-TextValue $ws.Range("D14") "0.09908"
Set-TextValue $ws.Range("E14") "-0.23%"
Set-TextValue $ws.Range("D15") "0.001432"
Set-TextValue $ws.Range("E15") "-1.18%"
Set-TextValue $ws.Range("D16") "0.005667"
Set-TextValue $ws.Range("E16") "-0.72%"
Set-TextValue $ws.Range("D17") "3.466"
Set-TextValue $ws.Range("E17") "-0.04%"
Set-TextValue $ws.Range("D18") "2.908"
Set-TextValue $ws.Range("E18") "10.47%"
Set-TextValue $ws.Range("D19") "0.3413"
Set-TextValue $ws.Range("E19") "0.67%"
Set-TextValue $ws.Range("D20") "0.1341"
Set-TextValue $ws.Range("E20") "-0.41%"
Set-TextValue $ws.Range("D21") "5.096"
Set-TextValue $ws.Range("E21") "-1.73%"
Set-TextValue $ws.Range("D22") "0.2215"
Set-TextValue $ws.Range("E22") "9.65%"
Set-TextValue $ws.Range("D23") "0.04575"
Set-TextValue $ws.Range("E23") "0.41%"
Set-TextValue $ws.Range("D24") "0.001248"
Set-TextValue $ws.Range("E24") "0.87%"
Set-TextValue $ws.Range("D25") "0.004698"
Set-TextValue $ws.Range("E25") "-6.73%"
Set-TextValue $ws.Range("D26") "0.0001252"
Set-TextValue $ws.Range("E26") "-21.91%"
Set-TextValue $ws.Range("D27") "0.0004509"
Set-TextValue $ws.Range("E27") "-5.09%"
Set-TextValue $ws.Range("D39") "0.01962"
Set-TextValue $ws.Range("E39") "5.49%"
Set-TextValue $ws.Range("D40") "0.04845"
Set-TextValue $ws.Range("E40") "1.63%"
Set-TextValue $ws.Range("D41") "0.007779"
Set-TextValue $ws.Range("E41") "2.16%"
Set-TextValue $ws.Range("D42") "0.1391"
Set-TextValue $ws.Range("E42") "-0.64%"
Set-TextValue $ws.Range("D43") "0.007852"
Set-TextValue $ws.Range("E43") "0.99%"
Set-TextValue $ws.Range("D44") "0.002113"
Set-TextValue $ws.Range("E44") "-4.18%"
Set-TextValue $ws.Range("D45") "0.01180"
Set-TextValue $ws.Range("E45") "6.62%"
Set-TextValue $ws.Range("D46") "0.00006661"
Set-TextValue $ws.Range("E46") "5.20%"
Set-TextValue $ws.Range("E47") "0.16%"
Set-TextValue $ws.Range("D48") "37.86"
Set-TextValue $ws.Range("E48") "-22.13%"
Set-TextValue $ws.Range("D49") "0.001703"
Set-TextValue $ws.Range("E49") "-14.85%"
Set-TextValue $ws.Range("D50") "0.00002104"
Set-TextValue $ws.Range("E50") "0.16%"
Set-TextValue $ws.Range("D51") "0.0002004"
Set-TextValue $ws.Range("E51") "0.16%"
